# Auto-generated edit script
# Applies updated market-price derived values (H,I,J,K,L,M,N columns)
# for specific rows across multiple worksheets, per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(82, 8).Value = 5372.8335   # H82: 3330 -> 5372.8335
$ws.Cells.Item(82, 9).Value = 118.5   # I82: 216.66667 -> 118.5
$ws.Cells.Item(82, 11).Value = 355.5   # K82: 650.00001 -> 355.5
$ws.Cells.Item(82, 13).Value = 50.5   # M82: -244.00001 -> 50.5
$ws.Cells.Item(85, 8).Value = 5372.8335   # H85: 3330 -> 5372.8335
$ws.Cells.Item(85, 9).Value = 118.5   # I85: 216.66667 -> 118.5
$ws.Cells.Item(85, 11).Value = 355.5   # K85: 650.00001 -> 355.5
$ws.Cells.Item(85, 13).Value = 1048.5   # M85: 753.99999 -> 1048.5
$ws.Cells.Item(112, 8).Value = 10871144   # H112: 15689282 -> 10871144
$ws.Cells.Item(112, 9).Value = 2466.6667   # I112: 2483.3333 -> 2466.6667
$ws.Cells.Item(112, 10).Value = 14707148   # J112: 25101360 -> 14707148
$ws.Cells.Item(112, 11).Value = 7400.000100000001   # K112: 7449.999899999999 -> 7400.000100000001
$ws.Cells.Item(112, 12).Value = 44121444   # L112: 75304080 -> 44121444
$ws.Cells.Item(112, 13).Value = -6292.000100000001   # M112: -6341.999899999999 -> -6292.000100000001
$ws.Cells.Item(112, 14).Value = -44123660   # N112: -75306296 -> -44123660
$ws.Cells.Item(113, 8).Value = 2929.0881   # H113: 3147.3872 -> 2929.0881
$ws.Cells.Item(113, 9).Value = 2956.8845   # I113: 3074.36 -> 2956.8845
$ws.Cells.Item(113, 10).Value = 2838.75   # J113: 3451.6667 -> 2838.75
$ws.Cells.Item(113, 11).Value = 2956.8845   # K113: 3074.36 -> 2956.8845
$ws.Cells.Item(113, 12).Value = 2838.75   # L113: 3451.6667 -> 2838.75
$ws.Cells.Item(113, 13).Value = 297.1154999999999   # M113: 179.6399999999999 -> 297.1154999999999
$ws.Cells.Item(113, 14).Value = -9346.75   # N113: -9959.6667 -> -9346.75
$ws.Cells.Item(132, 8).Value = 3849205.2   # H132: 2668844.8 -> 3849205.2
$ws.Cells.Item(132, 9).Value = 4881068   # I132: 3176636.2 -> 4881068
$ws.Cells.Item(132, 10).Value = 3171   # J132: 2940.5833 -> 3171
$ws.Cells.Item(132, 11).Value = 14643204   # K132: 9529908.600000001 -> 14643204
$ws.Cells.Item(132, 12).Value = 9513   # L132: 8821.749899999999 -> 9513
$ws.Cells.Item(132, 13).Value = -14640674   # M132: -9527378.600000001 -> -14640674
$ws.Cells.Item(132, 14).Value = -14573   # N132: -13881.7499 -> -14573
$ws.Cells.Item(137, 8).Value = 3382.3958   # H137: 2883.1404 -> 3382.3958
$ws.Cells.Item(137, 9).Value = 3566.6758   # I137: 2911.9783 -> 3566.6758
$ws.Cells.Item(137, 11).Value = 10700.0274   # K137: 8735.9349 -> 10700.0274
$ws.Cells.Item(137, 13).Value = -8150.027399999999   # M137: -6185.9349 -> -8150.027399999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1674.7407   # H61: 1136.766 -> 1674.7407
$ws.Cells.Item(61, 9).Value = 835.56525   # I61: 641.1429000000001 -> 835.56525
$ws.Cells.Item(61, 10).Value = 6500   # J61: 5300 -> 6500
$ws.Cells.Item(61, 11).Value = 835.56525   # K61: 641.1429000000001 -> 835.56525
$ws.Cells.Item(61, 12).Value = 6500   # L61: 5300 -> 6500
$ws.Cells.Item(61, 13).Value = -623.56525   # M61: -429.1429000000001 -> -623.56525
$ws.Cells.Item(61, 14).Value = -6924   # N61: -5724 -> -6924
$ws.Cells.Item(74, 8).Value = 841.2222   # H74: 974.5417 -> 841.2222
$ws.Cells.Item(74, 9).Value = 669.38464   # I74: 790.4091 -> 669.38464
$ws.Cells.Item(74, 10).Value = 1288   # J74: 3000 -> 1288
$ws.Cells.Item(74, 11).Value = 669.38464   # K74: 790.4091 -> 669.38464
$ws.Cells.Item(74, 12).Value = 1288   # L74: 3000 -> 1288
$ws.Cells.Item(74, 13).Value = 204.61536   # M74: 83.59090000000003 -> 204.61536
$ws.Cells.Item(74, 14).Value = -3036   # N74: -4748 -> -3036
$ws.Cells.Item(77, 8).Value = 841.2222   # H77: 974.5417 -> 841.2222
$ws.Cells.Item(77, 9).Value = 669.38464   # I77: 790.4091 -> 669.38464
$ws.Cells.Item(77, 10).Value = 1288   # J77: 3000 -> 1288
$ws.Cells.Item(77, 11).Value = 3346.9232   # K77: 3952.0455 -> 3346.9232
$ws.Cells.Item(77, 12).Value = 6440   # L77: 15000 -> 6440
$ws.Cells.Item(77, 13).Value = 1021.0768   # M77: 415.9545000000003 -> 1021.0768
$ws.Cells.Item(77, 14).Value = -15176   # N77: -23736 -> -15176
$ws.Cells.Item(132, 8).Value = 2736.7693   # H132: 2337.652 -> 2736.7693
$ws.Cells.Item(132, 9).Value = 2055.5186   # I132: 1655.7941 -> 2055.5186
$ws.Cells.Item(132, 11).Value = 6166.5558   # K132: 4967.3823 -> 6166.5558
$ws.Cells.Item(132, 13).Value = -3636.5558   # M132: -2437.3823 -> -3636.5558
$ws.Cells.Item(136, 8).Value = 1674.7407   # H136: 1136.766 -> 1674.7407
$ws.Cells.Item(136, 9).Value = 835.56525   # I136: 641.1429000000001 -> 835.56525
$ws.Cells.Item(136, 10).Value = 6500   # J136: 5300 -> 6500
$ws.Cells.Item(136, 11).Value = 2506.69575   # K136: 1923.4287 -> 2506.69575
$ws.Cells.Item(136, 12).Value = 19500   # L136: 15900 -> 19500
$ws.Cells.Item(136, 13).Value = 43.30425000000014   # M136: 626.5712999999998 -> 43.30425000000014
$ws.Cells.Item(136, 14).Value = -24600   # N136: -21000 -> -24600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2348.0981   # H107: 1923.1791 -> 2348.0981
$ws.Cells.Item(107, 9).Value = 1874.7297   # I107: 1579.6875 -> 1874.7297
$ws.Cells.Item(107, 10).Value = 3599.1428   # J107: 2790.9473 -> 3599.1428
$ws.Cells.Item(107, 11).Value = 1874.7297   # K107: 1579.6875 -> 1874.7297
$ws.Cells.Item(107, 12).Value = 3599.1428   # L107: 2790.9473 -> 3599.1428
$ws.Cells.Item(107, 13).Value = 45.27029999999991   # M107: 340.3125 -> 45.27029999999991
$ws.Cells.Item(107, 14).Value = -7439.1428   # N107: -6630.9473 -> -7439.1428
$ws.Cells.Item(134, 8).Value = 3337.6858   # H134: 3081.9211 -> 3337.6858
$ws.Cells.Item(134, 9).Value = 2996.96   # I134: 2686.3572 -> 2996.96
$ws.Cells.Item(134, 11).Value = 8990.880000000001   # K134: 8059.071599999999 -> 8990.880000000001
$ws.Cells.Item(134, 13).Value = -6455.880000000001   # M134: -5524.071599999999 -> -6455.880000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2396.4656   # H31: 2557.0942 -> 2396.4656
$ws.Cells.Item(31, 9).Value = 1402   # I31: 1518.7727 -> 1402
$ws.Cells.Item(31, 10).Value = 7169.9   # J31: 7633.3335 -> 7169.9
$ws.Cells.Item(31, 11).Value = 1402   # K31: 1518.7727 -> 1402
$ws.Cells.Item(31, 12).Value = 7169.9   # L31: 7633.3335 -> 7169.9
$ws.Cells.Item(31, 13).Value = -1107   # M31: -1223.7727 -> -1107
$ws.Cells.Item(31, 14).Value = -7759.9   # N31: -8223.333500000001 -> -7759.9
$ws.Cells.Item(34, 8).Value = 2396.4656   # H34: 2557.0942 -> 2396.4656
$ws.Cells.Item(34, 9).Value = 1402   # I34: 1518.7727 -> 1402
$ws.Cells.Item(34, 10).Value = 7169.9   # J34: 7633.3335 -> 7169.9
$ws.Cells.Item(34, 11).Value = 1402   # K34: 1518.7727 -> 1402
$ws.Cells.Item(34, 12).Value = 7169.9   # L34: 7633.3335 -> 7169.9
$ws.Cells.Item(34, 13).Value = -1200   # M34: -1316.7727 -> -1200
$ws.Cells.Item(34, 14).Value = -7573.9   # N34: -8037.3335 -> -7573.9
$ws.Cells.Item(58, 8).Value = 6946696   # H58: 7044629.5 -> 6946696
$ws.Cells.Item(58, 9).Value = 1407.9445   # I58: 1473.0385 -> 1407.9445
$ws.Cells.Item(58, 10).Value = 27782560   # J58: 26320636 -> 27782560
$ws.Cells.Item(58, 11).Value = 1407.9445   # K58: 1473.0385 -> 1407.9445
$ws.Cells.Item(58, 12).Value = 27782560   # L58: 26320636 -> 27782560
$ws.Cells.Item(58, 13).Value = -1204.9445   # M58: -1270.0385 -> -1204.9445
$ws.Cells.Item(58, 14).Value = -27782966   # N58: -26321042 -> -27782966
$ws.Cells.Item(99, 8).Value = 1702.9524   # H99: 1702.1818 -> 1702.9524
$ws.Cells.Item(99, 9).Value = 1118.25   # I99: 1201.6471 -> 1118.25
$ws.Cells.Item(99, 10).Value = 3574   # J99: 3404 -> 3574
$ws.Cells.Item(99, 11).Value = 1118.25   # K99: 1201.6471 -> 1118.25
$ws.Cells.Item(99, 12).Value = 3574   # L99: 3404 -> 3574
$ws.Cells.Item(99, 13).Value = 379.75   # M99: 296.3529000000001 -> 379.75
$ws.Cells.Item(99, 14).Value = -6570   # N99: -6400 -> -6570
$ws.Cells.Item(126, 8).Value = 1702.9524   # H126: 1702.1818 -> 1702.9524
$ws.Cells.Item(126, 9).Value = 1118.25   # I126: 1201.6471 -> 1118.25
$ws.Cells.Item(126, 10).Value = 3574   # J126: 3404 -> 3574
$ws.Cells.Item(126, 11).Value = 3354.75   # K126: 3604.9413 -> 3354.75
$ws.Cells.Item(126, 12).Value = 10722   # L126: 10212 -> 10722
$ws.Cells.Item(126, 13).Value = -884.75   # M126: -1134.9413 -> -884.75
$ws.Cells.Item(126, 14).Value = -15662   # N126: -15152 -> -15662
$ws.Cells.Item(132, 8).Value = 2150.7273   # H132: 2040.2979 -> 2150.7273
$ws.Cells.Item(132, 9).Value = 1670.875   # I132: 1553.6111 -> 1670.875
$ws.Cells.Item(132, 10).Value = 3430.3333   # J132: 3633.0908 -> 3430.3333
$ws.Cells.Item(132, 11).Value = 5012.625   # K132: 4660.8333 -> 5012.625
$ws.Cells.Item(132, 12).Value = 10290.9999   # L132: 10899.2724 -> 10290.9999
$ws.Cells.Item(132, 13).Value = -2482.625   # M132: -2130.8333 -> -2482.625
$ws.Cells.Item(132, 14).Value = -15350.9999   # N132: -15959.2724 -> -15350.9999
$ws.Cells.Item(134, 8).Value = 2093.9312   # H134: 1579.262 -> 2093.9312
$ws.Cells.Item(134, 9).Value = 1192.8096   # I134: 878.6857 -> 1192.8096
$ws.Cells.Item(134, 10).Value = 4459.375   # J134: 5082.143 -> 4459.375
$ws.Cells.Item(134, 11).Value = 3578.4288   # K134: 2636.0571 -> 3578.4288
$ws.Cells.Item(134, 12).Value = 13378.125   # L134: 15246.429 -> 13378.125
$ws.Cells.Item(134, 13).Value = -1043.4288   # M134: -101.0571 -> -1043.4288
$ws.Cells.Item(134, 14).Value = -18448.125   # N134: -20316.429 -> -18448.125
$ws.Cells.Item(136, 8).Value = 6946696   # H136: 7044629.5 -> 6946696
$ws.Cells.Item(136, 9).Value = 1407.9445   # I136: 1473.0385 -> 1407.9445
$ws.Cells.Item(136, 10).Value = 27782560   # J136: 26320636 -> 27782560
$ws.Cells.Item(136, 11).Value = 4223.833500000001   # K136: 4419.1155 -> 4223.833500000001
$ws.Cells.Item(136, 12).Value = 83347680   # L136: 78961908 -> 83347680
$ws.Cells.Item(136, 13).Value = -1673.833500000001   # M136: -1869.1155 -> -1673.833500000001
$ws.Cells.Item(136, 14).Value = -83352780   # N136: -78967008 -> -83352780

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 80436.57000000001   # H69: 75094.13 -> 80436.57000000001
$ws.Cells.Item(69, 9).Value = 906   # I69: 704 -> 906
$ws.Cells.Item(69, 11).Value = 2718   # K69: 2112 -> 2718
$ws.Cells.Item(69, 13).Value = -1907   # M69: -1301 -> -1907
$ws.Cells.Item(72, 8).Value = 80436.57000000001   # H72: 75094.13 -> 80436.57000000001
$ws.Cells.Item(72, 9).Value = 906   # I72: 704 -> 906
$ws.Cells.Item(72, 11).Value = 8154   # K72: 6336 -> 8154
$ws.Cells.Item(72, 13).Value = -4098   # M72: -2280 -> -4098
$ws.Cells.Item(123, 8).Value = 2990   # H123: 3266.3635 -> 2990
$ws.Cells.Item(123, 9).Value = 340   # I123: 465 -> 340
$ws.Cells.Item(123, 10).Value = 3520   # J123: 3888.889 -> 3520
$ws.Cells.Item(123, 11).Value = 1020   # K123: 1395 -> 1020
$ws.Cells.Item(123, 12).Value = 10560   # L123: 11666.667 -> 10560
$ws.Cells.Item(123, 13).Value = 1430   # M123: 1055 -> 1430
$ws.Cells.Item(123, 14).Value = -15460   # N123: -16566.667 -> -15460
$ws.Cells.Item(129, 8).Value = 23265   # H129: 28934.316 -> 23265
$ws.Cells.Item(129, 9).Value = 3150.9092   # I129: 4118.5713 -> 3150.9092
$ws.Cells.Item(129, 10).Value = 40284.617   # J129: 43410.168 -> 40284.617
$ws.Cells.Item(129, 11).Value = 9452.7276   # K129: 12355.7139 -> 9452.7276
$ws.Cells.Item(129, 12).Value = 120853.851   # L129: 130230.504 -> 120853.851
$ws.Cells.Item(129, 13).Value = -4452.7276   # M129: -7355.713899999999 -> -4452.7276
$ws.Cells.Item(129, 14).Value = -130853.851   # N129: -140230.504 -> -130853.851
$ws.Cells.Item(130, 8).Value = 2600   # H130: 3000 -> 2600
$ws.Cells.Item(130, 10).Value = 2600   # J130: 3000 -> 2600
$ws.Cells.Item(130, 12).Value = 7800   # L130: 9000 -> 7800
$ws.Cells.Item(130, 14).Value = -17840   # N130: -19040 -> -17840
$ws.Cells.Item(131, 8).Value = 1923.4062   # H131: 2184.76 -> 1923.4062
$ws.Cells.Item(131, 10).Value = 1686.1428   # J131: 2034.2142 -> 1686.1428
$ws.Cells.Item(131, 12).Value = 5058.428400000001   # L131: 6102.642599999999 -> 5058.428400000001
$ws.Cells.Item(131, 14).Value = -15138.4284   # N131: -16182.6426 -> -15138.4284
$ws.Cells.Item(133, 8).Value = 4043.842   # H133: 4526.4375 -> 4043.842
$ws.Cells.Item(133, 9).Value = 5233.75   # I133: 5827.143 -> 5233.75
$ws.Cells.Item(133, 10).Value = 3178.4546   # J133: 3514.7778 -> 3178.4546
$ws.Cells.Item(133, 11).Value = 15701.25   # K133: 17481.429 -> 15701.25
$ws.Cells.Item(133, 12).Value = 9535.363799999999   # L133: 10544.3334 -> 9535.363799999999
$ws.Cells.Item(133, 13).Value = -10641.25   # M133: -12421.429 -> -10641.25
$ws.Cells.Item(133, 14).Value = -19655.3638   # N133: -20664.3334 -> -19655.3638
$ws.Cells.Item(134, 8).Value = 3384.3845   # H134: 2452 -> 3384.3845
$ws.Cells.Item(134, 9).Value = 1572.5   # I134: 1595.3334 -> 1572.5
$ws.Cells.Item(134, 10).Value = 4937.4287   # J134: 3994 -> 4937.4287
$ws.Cells.Item(134, 11).Value = 4717.5   # K134: 4786.0002 -> 4717.5
$ws.Cells.Item(134, 12).Value = 14812.2861   # L134: 11982 -> 14812.2861
$ws.Cells.Item(134, 13).Value = 352.5   # M134: 283.9997999999996 -> 352.5
$ws.Cells.Item(134, 14).Value = -24952.2861   # N134: -22122 -> -24952.2861
$ws.Cells.Item(136, 8).Value = 1757.5   # H136: 1687.6 -> 1757.5
$ws.Cells.Item(136, 10).Value = 2662.2222   # J136: 2397 -> 2662.2222
$ws.Cells.Item(136, 12).Value = 7986.6666   # L136: 7191 -> 7986.6666
$ws.Cells.Item(136, 14).Value = -18186.6666   # N136: -17391 -> -18186.6666
$ws.Cells.Item(138, 8).Value = 2660.0833   # H138: 3068.5334 -> 2660.0833
$ws.Cells.Item(138, 9).Value = 1058.5   # I138: 1247.2858 -> 1058.5
$ws.Cells.Item(138, 10).Value = 4261.6665   # J138: 4662.125 -> 4261.6665
$ws.Cells.Item(138, 11).Value = 3175.5   # K138: 3741.8574 -> 3175.5
$ws.Cells.Item(138, 12).Value = 12784.9995   # L138: 13986.375 -> 12784.9995
$ws.Cells.Item(138, 13).Value = 1964.5   # M138: 1398.1426 -> 1964.5
$ws.Cells.Item(138, 14).Value = -23064.9995   # N138: -24266.375 -> -23064.9995
$ws.Cells.Item(139, 8).Value = 13893092   # H139: 12503684 -> 13893092
$ws.Cells.Item(139, 9).Value = 22731070   # I139: 27781426 -> 22731070
$ws.Cells.Item(139, 10).Value = 4841.4287   # J139: 3713.2727 -> 4841.4287
$ws.Cells.Item(139, 11).Value = 68193210   # K139: 83344278 -> 68193210
$ws.Cells.Item(139, 12).Value = 14524.2861   # L139: 11139.8181 -> 14524.2861
$ws.Cells.Item(139, 13).Value = -68188070   # M139: -83339138 -> -68188070
$ws.Cells.Item(139, 14).Value = -24804.2861   # N139: -21419.8181 -> -24804.2861
$ws.Cells.Item(140, 8).Value = 16671636   # H140: 23813194 -> 16671636
$ws.Cells.Item(140, 10).Value = 6828.5713   # J140: 5950 -> 6828.5713
$ws.Cells.Item(140, 12).Value = 20485.7139   # L140: 17850 -> 20485.7139
$ws.Cells.Item(140, 14).Value = -30845.7139   # N140: -28210 -> -30845.7139
$ws.Cells.Item(141, 8).Value = 3950   # H141: 4000 -> 3950
$ws.Cells.Item(141, 9).Value = 3950   # I141: 4000 -> 3950
$ws.Cells.Item(141, 11).Value = 11850   # K141: 12000 -> 11850
$ws.Cells.Item(141, 13).Value = -6670   # M141: -6820 -> -6670

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4723.161   # H132: 4528.724 -> 4723.161
$ws.Cells.Item(132, 9).Value = 5496.1763   # I132: 5410.067 -> 5496.1763
$ws.Cells.Item(132, 10).Value = 3784.5   # J132: 3584.4285 -> 3784.5
$ws.Cells.Item(132, 11).Value = 16488.5289   # K132: 16230.201 -> 16488.5289
$ws.Cells.Item(132, 12).Value = 11353.5   # L132: 10753.2855 -> 11353.5
$ws.Cells.Item(132, 13).Value = -13958.5289   # M132: -13700.201 -> -13958.5289
$ws.Cells.Item(132, 14).Value = -16413.5   # N132: -15813.2855 -> -16413.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1977.4615   # H7: 2059 -> 1977.4615
$ws.Cells.Item(7, 9).Value = 1125.25   # I7: 1143.2858 -> 1125.25
$ws.Cells.Item(7, 11).Value = 1125.25   # K7: 1143.2858 -> 1125.25
$ws.Cells.Item(7, 13).Value = -1013.25   # M7: -1031.2858 -> -1013.25
$ws.Cells.Item(104, 8).Value = 27952.715   # H104: 30000 -> 27952.715
$ws.Cells.Item(104, 10).Value = 27952.715   # J104: 30000 -> 27952.715
$ws.Cells.Item(104, 12).Value = 27952.715   # L104: 30000 -> 27952.715
$ws.Cells.Item(104, 14).Value = -34940.715   # N104: -36988 -> -34940.715
$ws.Cells.Item(124, 8).Value = 40000   # H124: 36666.668 -> 40000
$ws.Cells.Item(124, 10).Value = 40000   # J124: 36666.668 -> 40000
$ws.Cells.Item(124, 12).Value = 40000   # L124: 36666.668 -> 40000
$ws.Cells.Item(124, 14).Value = -49820   # N124: -46486.668 -> -49820
$ws.Cells.Item(126, 8).Value = 1977.4615   # H126: 2059 -> 1977.4615
$ws.Cells.Item(126, 9).Value = 1125.25   # I126: 1143.2858 -> 1125.25
$ws.Cells.Item(126, 11).Value = 3375.75   # K126: 3429.8574 -> 3375.75
$ws.Cells.Item(126, 13).Value = -905.75   # M126: -959.8574000000003 -> -905.75
$ws.Cells.Item(132, 8).Value = 1668.7572   # H132: 1807 -> 1668.7572
$ws.Cells.Item(132, 9).Value = 1054.037   # I132: 1146.9375 -> 1054.037
$ws.Cells.Item(132, 10).Value = 3743.4375   # J132: 3787.1875 -> 3743.4375
$ws.Cells.Item(132, 11).Value = 3162.111   # K132: 3440.8125 -> 3162.111
$ws.Cells.Item(132, 12).Value = 11230.3125   # L132: 11361.5625 -> 11230.3125
$ws.Cells.Item(132, 13).Value = -632.1109999999999   # M132: -910.8125 -> -632.1109999999999
$ws.Cells.Item(132, 14).Value = -16290.3125   # N132: -16421.5625 -> -16290.3125
$ws.Cells.Item(136, 8).Value = 1645.3334   # H136: 1466.6492 -> 1645.3334
$ws.Cells.Item(136, 9).Value = 1134.3864   # I136: 1027.18 -> 1134.3864
$ws.Cells.Item(136, 10).Value = 4857   # J136: 4605.7144 -> 4857
$ws.Cells.Item(136, 11).Value = 3403.1592   # K136: 3081.54 -> 3403.1592
$ws.Cells.Item(136, 12).Value = 14571   # L136: 13817.1432 -> 14571
$ws.Cells.Item(136, 13).Value = -853.1592000000001   # M136: -531.54 -> -853.1592000000001
$ws.Cells.Item(136, 14).Value = -19671   # N136: -18917.1432 -> -19671

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 7685.1646   # H132: 8162.175 -> 7685.1646
$ws.Cells.Item(132, 9).Value = 1784.4906   # I132: 1862.0197 -> 1784.4906
$ws.Cells.Item(132, 10).Value = 17458.156   # J132: 19241.758 -> 17458.156
$ws.Cells.Item(132, 11).Value = 5353.4718   # K132: 5586.0591 -> 5353.4718
$ws.Cells.Item(132, 12).Value = 52374.46799999999   # L132: 57725.274 -> 52374.46799999999
$ws.Cells.Item(132, 13).Value = -2823.4718   # M132: -3056.0591 -> -2823.4718
$ws.Cells.Item(132, 14).Value = -57434.46799999999   # N132: -62785.274 -> -57434.46799999999
$ws.Cells.Item(136, 8).Value = 839.8570999999999   # H136: 804.625 -> 839.8570999999999
$ws.Cells.Item(136, 9).Value = 627.0833   # I136: 542.7586 -> 627.0833
$ws.Cells.Item(136, 10).Value = 1304.091   # J136: 1495 -> 1304.091
$ws.Cells.Item(136, 11).Value = 1881.2499   # K136: 1628.2758 -> 1881.2499
$ws.Cells.Item(136, 12).Value = 3912.273   # L136: 4485 -> 3912.273
$ws.Cells.Item(136, 13).Value = 668.7501   # M136: 921.7242000000001 -> 668.7501
$ws.Cells.Item(136, 14).Value = -9012.272999999999   # N136: -9585 -> -9012.272999999999

Write-Host "Applied 253 cell updates across 8 sheets."